$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.217.34"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.894.12"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "'307.07"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").Value = "'0.5221"
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("D8").Value = "'0.3752"
$ws.Range("E8").Value = "  -1.01%  "
$ws.Range("D9").Value = "'0.07262"
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("D10").Value = "'21.18"
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("D11").Value = "'0.8981"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "'0.08180"
$ws.Range("E12").Value = "  +6.63%  "
$ws.Range("D13").Value = "'96.80"
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("D14").Value = "1.892.48"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").Value = "'5.273"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "'0.000008590"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "'1.002"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").Value = "27.241.66"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").Value = "'5.080"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").Value = "'6.402"
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("D24").Value = "'147.54"
$ws.Range("E24").Value = "  +1.10%  "
$ws.Range("D25").Value = "'2.286"
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'18.19"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "'1.743"
$ws.Range("E27").Value = "  +0.74%  "
$ws.Range("D28").Value = "'114.95"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("D29").Value = "'4.918"
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("D31").Value = "'0.09228"
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("D32").Value = "'0.05040"
$ws.Range("E32").Value = "  -0.84%  "
$ws.Range("D33").Value = "'0.7906"
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("D34").Value = "'1.216"
$ws.Range("E34").Value = "  -2.63%  "
$ws.Range("D35").Value = "'3.435"
$ws.Range("E35").Value = "  +4.07%  "
$ws.Range("D36").Value = "'2.960"
$ws.Range("E36").Value = "  -1.06%  "
$ws.Range("D37").Value = "'2.566"
$ws.Range("E37").Value = "  -2.13%  "
$ws.Range("D38").Value = "'0.5651"
$ws.Range("E38").Value = "  -0.52%  "
$ws.Range("D39").Value = "'0.01983"
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("D40").Value = "'1.074"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").Value = "'8.925"
$ws.Range("E41").Value = "  -1.09%  "
$ws.Range("D42").Value = "'6.533"
$ws.Range("E42").Value = "  -1.60%  "
$ws.Range("D43").Value = "'115.22"
$ws.Range("E43").Value = "  -3.18%  "
$ws.Range("D44").Value = "'0.1515"
$ws.Range("E44").Value = "  -0.51%  "
$ws.Range("D45").Value = "'0.4859"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "'1.002"
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("D47").Value = "'10.05"
$ws.Range("E47").Value = "  -1.56%  "
$ws.Range("D48").Value = "'1.616"
$ws.Range("E48").Value = "  +0.57%  "
$ws.Range("D49").Value = "'38.11"
$ws.Range("E49").Value = "  +1.77%  "
$ws.Range("D50").Value = "'63.28"
$ws.Range("E50").Value = "  -1.69%  "
$ws.Range("D51").Value = "'0.05941"
$ws.Range("E51").Value = "  +0.26%  "
